# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 (the financial-documents table) switches from
#    the deck's custom table style to the built-in
#    "{1159D0E1-D28F-4D22-9C2B-9D8795FE03EA}" table style.
#
# 2) The theme that actually paints the deck (slideMaster1 -> theme2.xml,
#    which is also the presentation's primary theme relationship) swaps
#    its 12-colour theme palette from the "Red Violet"/Integral colours
#    it has now to the plain "Office" colours that the (otherwise unused,
#    notes-master-only) theme1.xml currently carries. Font scheme and
#    format scheme are already identical between the two theme parts, so
#    only the colour scheme actually needs to move.

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
    }
}
$table = $tableShape.Table
$table.ApplyStyle("{1159D0E1-D28F-4D22-9C2B-9D8795FE03EA}")

# --- 2) Theme colours -------------------------------------------------
# Order exposed by ThemeColorScheme is: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink - the same order as <a:clrScheme> in the OOXML.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
